$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q4" right after "总计" (i.e. before the
#    current "2022-Q3" sheet), pushing every later sheet one tab to the right
#    (Excel keeps each existing sheet's name+content attached, just its tab
#    position moves).
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# Header row (same layout as the other quarterly sheets).
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

# Row 2 - 景顺长城中证500行业中性低波动指数
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").Value = "'003318"
$q4Sheet.Range("C2").Value = "景顺长城中证500行业中性低波动指数"
$q4Sheet.Range("D2").Value = "'10.01"
$q4Sheet.Range("E2").Value = "'93.81"
$q4Sheet.Range("F2").Value = "'1.05"
$q4Sheet.Range("G2").Value = "'0.1051"
$q4Sheet.Range("H2").Value = 7

# Row 3 - 华安中证500行业中性低波动ETF
$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").Value = "'512260"
$q4Sheet.Range("C3").Value = "华安中证500行业中性低波动ETF"
$q4Sheet.Range("D3").Value = "'0.94"
$q4Sheet.Range("E3").Value = "'97.66"
$q4Sheet.Range("F3").Value = "'1.10"
$q4Sheet.Range("G3").Value = "'0.0103"
$q4Sheet.Range("H3").Value = 7

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new row 2 for 2022-Q4 and
#    shift the existing quarter rows down by one.
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.12

# Re-number the index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

# Line up the new row's look with the rest of the table: column A keeps the
# bold/centred "index" style, columns B:D stay unstyled (same as every other
# data row).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
